$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) New header cell F1 = "Files"
# ---------------------------------------------------------------
$ws.Range("F1").Value = "Files"

# ---------------------------------------------------------------
# 2) Column D gets wider (to fit wrapped "Files" links)
# ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 26.45

# ---------------------------------------------------------------
# 3) New row 7 - start by copying row 6 so borders/fills/fonts
#    (cell styles) for A/B/C/E line up with the rest of the table.
# ---------------------------------------------------------------
$ws.Range("A6:F6").Copy($ws.Range("A7:F7"))
$ws.Rows.Item(7).RowHeight = 46.65

$ws.Range("A7").Value = "Ruby Necklace"
$ws.Range("B7").Value = "Natural (33 pieces) some are Red & some are Pigeon Blood color"
$ws.Range("C7").Value = "USD - 4 millions"

# ---------------------------------------------------------------
# D7: single bracketed link, wraps text
# ---------------------------------------------------------------
$dLink = "https://drive.google.com/file/d/1wBiiQ9sF0-rWSixLS3qtCamImi_FIiKi/view?usp=drive_link"
$dTxt = "[" + $dLink + "]"
$ws.Range("D7").Value = $dTxt
$ws.Range("D7").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("D7"), $dLink, "", "", $dTxt) | Out-Null

$dLen = $dTxt.Length
$dPart1 = $ws.Range("D7").Characters(1, $dLen - 1)
$dPart1.Font.Name = "Arial"
$dPart1.Font.Size = 10
$dPart1.Font.Underline = $true
$dPart1.Font.ColorIndex = 11
$dPart2 = $ws.Range("D7").Characters($dLen, 1)
$dPart2.Font.Name = "Arial"
$dPart2.Font.Size = 10
$dPart2.Font.Underline = $true
$dPart2.Font.ColorIndex = 11

# ---------------------------------------------------------------
# E7: "[" + link + "]" where only the link itself is hyperlinked /
#     colored; the brackets stay plain black text.
# ---------------------------------------------------------------
$eLink = "https://drive.google.com/file/d/1k9fM7cg01YITiitPDL6Jda50TLNYT0NO/view?usp=drive_link"
$eTxt = "[" + $eLink + "]"
$ws.Range("E7").Value = $eTxt
$ws.Hyperlinks.Add($ws.Range("E7"), $eLink, "", "", $eTxt) | Out-Null

$eP1 = $ws.Range("E7").Characters(1, 1)
$eP1.Font.Name = "Arial"
$eP1.Font.Size = 10
$eP1.Font.ColorIndex = 1
$eP1.Font.Underline = $false

$eP2 = $ws.Range("E7").Characters(2, $eLink.Length)
$eP2.Font.Name = "Arial"
$eP2.Font.Size = 10
$eP2.Font.Underline = $true
$eP2.Font.ColorIndex = 11

$eP3 = $ws.Range("E7").Characters(2 + $eLink.Length, 1)
$eP3.Font.Name = "Arial"
$eP3.Font.Size = 10
$eP3.Font.ColorIndex = 1
$eP3.Font.Underline = $false

# ---------------------------------------------------------------
# F7: single bracketed link (same pattern as D7) but without wrap -
#     reuse E4's cell style (numFmt text + border, no wrap).
# ---------------------------------------------------------------
$fLink = "https://drive.google.com/file/d/1hl6d9iMoFpet0_Y8u9n6TG5qcMScUp57/view?usp=drive_link"
$fTxt = "[" + $fLink + "]"
$ws.Range("E4").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("F7").Value = $fTxt
$ws.Hyperlinks.Add($ws.Range("F7"), $fLink, "", "", $fTxt) | Out-Null

$fLen = $fTxt.Length
$fPart1 = $ws.Range("F7").Characters(1, $fLen - 1)
$fPart1.Font.Name = "Arial"
$fPart1.Font.Size = 10
$fPart1.Font.Underline = $true
$fPart1.Font.ColorIndex = 11
$fPart2 = $ws.Range("F7").Characters($fLen, 1)
$fPart2.Font.Name = "Arial"
$fPart2.Font.Size = 10
$fPart2.Font.Underline = $true
$fPart2.Font.ColorIndex = 11

Write-Host "done"
